$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.506.33"
$ws.Range("E2").Value = "  +1.78%  "

# Row 3
$ws.Range("D3").Value = "3.258.44"
$ws.Range("E3").Value = "  -0.77%  "

# Row 4
$c = $ws.Range("D4")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$c = $ws.Range("D5")
$c.Value = "'566.14"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "

# Row 6
$c = $ws.Range("D6")
$c.Value = "'174.21"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -4.20%  "

# Row 7
$c = $ws.Range("D7")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$c = $ws.Range("D8")
$c.Value = "'0.580"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.86%  "

# Row 9
$ws.Range("D9").Value = "3.250.22"
$ws.Range("E9").Value = "  -0.88%  "

# Row 10
$ws.Range("E10").Value = "  -1.40%  "

# Row 11
$c = $ws.Range("D11")
$c.Value = "'0.566"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.80%  "

# Row 12
$c = $ws.Range("D12")
$c.Value = "'44.95"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.84%  "

# Row 13
$c = $ws.Range("D13")
$c.Value = "'0.0000266"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.63%  "

# Row 14
$c = $ws.Range("D14")
$c.Value = "'689.77"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +10.43%  "

# Row 15
$ws.Range("D15").Value = "3.778.03"
$ws.Range("E15").Value = "  -0.87%  "

# Row 16
$c = $ws.Range("D16")
$c.Value = "'8.24"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.45%  "

# Row 17
$ws.Range("D17").Value = "66.617.60"
$ws.Range("E17").Value = "  +1.73%  "

# Row 18
$ws.Range("E18").Value = "  +0.89%  "

# Row 19
$ws.Range("D19").Value = "3.259.05"
$ws.Range("E19").Value = "  -0.89%  "

# Row 20
$c = $ws.Range("D20")
$c.Value = "'17.18"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.28%  "

# Row 21
$c = $ws.Range("D21")
$c.Value = "'10.64"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.22%  "

# Row 22
$c = $ws.Range("D22")
$c.Value = "'0.880"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "

# Row 23
$c = $ws.Range("D23")
$c.Value = "'16.76"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -6.52%  "

# Row 24
$c = $ws.Range("D24")
$c.Value = "'5.08"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.78%  "

# Row 25
$c = $ws.Range("D25")
$c.Value = "'97.17"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.35%  "

# Row 26
$c = $ws.Range("D26")
$c.Value = "'3.84"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.61%  "

# Row 27
$c = $ws.Range("D27")
$c.Value = "'2.68"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "

# Row 28
$c = $ws.Range("D28")
$c.Value = "'9.19"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.28%  "

# Row 29
$c = $ws.Range("D29")
$c.Value = "'32.54"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.63%  "

# Row 30
$c = $ws.Range("D30")
$c.Value = "'8.32"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "

# Row 31
$c = $ws.Range("D31")
$c.Value = "'6.65"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.20%  "

# Row 32
$c = $ws.Range("D32")
$c.Value = "'572.93"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.04%  "

# Row 33
$ws.Range("D33").Value = "3.844.59"
$ws.Range("E33").Value = "  +0.33%  "

# Row 34
$c = $ws.Range("D34")
$c.Value = "'10.69"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D35")
$c.Value = "'0.102"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.84%  "

# Row 36
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D36")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$c = $ws.Range("D37")
$c.Value = "'55.15"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "

# Row 38
$c = $ws.Range("D38")
$c.Value = "'3.24"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -11.70%  "

# Row 39
$c = $ws.Range("D39")
$c.Value = "'0.128"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.05%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D40")
$c.Value = "'2.58"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.15%  "

# Row 41
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Range("D41")
$c.Value = "'3.33"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.51%  "

# Row 42
$c = $ws.Range("D42")
$c.Value = "'31.45"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "

# Row 43
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0660"
$ws.Range("E43").Value = "  -2.27%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D44")
$c.Value = "'3.02"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.67%  "

# Row 45
$ws.Range("E45").Value = "  -1.69%  "

# Row 46
$c = $ws.Range("D46")
$c.Value = "'0.0402"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "

# Row 47
$ws.Range("E47").Value = "  +0.29%  "

# Row 48
$c = $ws.Range("D48")
$c.Value = "'0.127"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "

# Row 49
$c = $ws.Range("D49")
$c.Value = "'1.35"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +7.93%  "

# Row 50
$c = $ws.Range("D50")
$c.Value = "'2.49"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "

# Row 51
$c = $ws.Range("D51")
$c.Value = "'129.03"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
